# Apply weekly update to Fruta/hortaliza data:
# Rotate the data among rows 2, 4 and 5:
#   Row 2 <- old Row 4
#   Row 4 <- old Row 5
#   Row 5 <- old Row 2
# (only columns D, J, K, L, M, P carry data that actually changes
#  per the cyclic rotation; other columns in these rows are identical
#  across the three rows already)

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$cols = @("D", "J", "K", "L", "M", "P")

# Capture original values before overwriting anything
$row2 = @{}
$row4 = @{}
$row5 = @{}
foreach ($col in $cols) {
    $row2[$col] = $ws.Range("${col}2").Value2
    $row4[$col] = $ws.Range("${col}4").Value2
    $row5[$col] = $ws.Range("${col}5").Value2
}

foreach ($col in $cols) {
    $ws.Range("${col}2").Value2 = $row4[$col]
    $ws.Range("${col}4").Value2 = $row5[$col]
    $ws.Range("${col}5").Value2 = $row2[$col]
}
